# Generate Report for Handoff
# Update "Latest Handoff Datetime" (column D) for the row corresponding to
# "68c3c0c3-8147-400f-94a1-fe4fb87aeb7b.md" (row 5) on both language sheets,
# reflecting that a new handoff was generated.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-03-09 20:35:51"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-03-09 20:35:56"
